{"js": "// 1. \"Note that signups are not first-come-first-serve.\" ->\n//    \"Note that signups are not done on a first-come, first-served basis.\"\n{\n  const results = context.document.body.search(\" first-come-first-serve.\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\" done on a first-come, first-served basis.\", \"Replace\");\n  }\n}\n\n// 2. \" conflict with this time, and \" -> \" conflict with this time and \" (drop the comma)\n{\n  const results = context.document.body.search(\" conflict with this time, and \", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\" conflict with this time and \", \"Replace\");\n  }\n}\n\n// 3. \"send a request by email to the Head TA\" -> \"fill out the Head TA's exam form\"\n{\n  const results = context.document.body.search(\"send a request by email to the Head TA\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"fill out the Head TA\\u2019s exam form\", \"Replace\");\n  }\n}\n\n// 4. Final exam date: \"August 18th\" -> \"August 17th\" (only the day number changes)\n{\n  const results = context.document.body.search(\"18th from\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"17th from\", \"Replace\");\n  }\n}\n\n// 5. Copyright line: remove the comma before \"and based on similar handouts\"\n{\n  const results = context.document.body.search(\", and based on similar handouts\", { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\" and based on similar handouts\", \"Replace\");\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$wdFindContinue = 1   # wdFindContinue\n$wdReplaceAll = 2     # wdReplaceAll\n$rightSingleQuote = [char]8217\n\n# 1. \"Note that signups are not first-come-first-serve.\" ->\n#    \"Note that signups are not done on a first-come, first-served basis.\"\n$rng = $d.Content\n$rng.Find.Execute(\n    \" first-come-first-serve.\",\n    $false, $false, $false, $false, $false,\n    $true, $wdFindContinue, $false,\n    \" done on a first-come, first-served basis.\",\n    $wdReplaceAll\n) | Out-Null\n\n# Word re-anchors its internal \"_GoBack\" (last-edit-location) bookmark to wherever\n# text was most recently edited; move it along with the edit above.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n$rng = $d.Content\n$found = $rng.Find.Execute(\"first-serve\", $false, $false, $false, $false, $false, $true, $wdFindContinue)\nif ($found) {\n    $bmRange = $rng.Duplicate\n    $bmRange.Collapse(1) | Out-Null  # wdCollapseStart\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange) | Out-Null\n}\n\n# 2. \" conflict with this time, and \" -> \" conflict with this time and \" (drop the comma)\n$rng = $d.Content\n$rng.Find.Execute(\n    \" conflict with this time, and \",\n    $false, $false, $false, $false, $false,\n    $true, $wdFindContinue, $false,\n    \" conflict with this time and \",\n    $wdReplaceAll\n) | Out-Null\n\n# 3. \"send a request by email to the Head TA\" -> \"fill out the Head TA's exam form\"\n$rng = $d.Content\n$rng.Find.Execute(\n    \"send a request by email to the Head TA\",\n    $false, $false, $false, $false, $false,\n    $true, $wdFindContinue, $false,\n    \"fill out the Head TA\" + $rightSingleQuote + \"s exam form\",\n    $wdReplaceAll\n) | Out-Null\n\n# 4. Final exam date: \"August 18th\" -> \"August 17th\" (only the day number changes)\n$rng = $d.Content\n$rng.Find.Execute(\n    \"18th from\",\n    $false, $false, $false, $false, $false,\n    $true, $wdFindContinue, $false,\n    \"17th from\",\n    $wdReplaceAll\n) | Out-Null\n\n# 5. Copyright line: remove the comma before \"and based on similar handouts\"\n$rng = $d.Content\n$rng.Find.Execute(\n    \", and based on similar handouts\",\n    $false, $false, $false, $false, $false,\n    $true, $wdFindContinue, $false,\n    \" and based on similar handouts\",\n    $wdReplaceAll\n) | Out-Null\n"}
